# Re-order the header columns in row 2 of the "PeopleReachedValues" staging
# template. The ID column moves to the front (right after the
# "For internal use only" notice in A1) and "Notes" moves to the end, with
# the rest of the BusinessKey columns kept in alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "PeopleReachedValuesID",
    "ActivityBusinessKey",
    "AgeBandBusinessKey",
    "CommunityTypeBusinessKey",
    "DataVersionBusinessKey",
    "DonorBusinessKey",
    "Framework_IndicatorBusinessKey",
    "FrameworkBusinessKey",
    "GenderBusinessKey",
    "GroupBusinessKey",
    "InstitutionBusinessKey",
    "LocationBusinessKey",
    "OutcomeBusinessKey",
    "OutputBusinessKey",
    "ProgrammeBusinessKey",
    "ProjectBusinessKey",
    "ReportingPeriodBusinessKey",
    "ResultAreaBusinessKey",
    "StatusTypeBusinessKey",
    "StrategicElementBusinessKey",
    "SubOutputBusinessKey",
    "Notes"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newHeaders[$i]
}
